$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3: columns J and O are removed in the updated data (error
# correction dropped these metrics), clear them individually.
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 2: corrected financial figures
$ws.Range("D2").Value = 2502
$ws.Range("E2").Value = 92
$ws.Range("F2").Value = 92
$ws.Range("G2").Value = 84
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 70
$ws.Range("K2").Value = 2352
$ws.Range("L2").Value = 873
$ws.Range("M2").Value = 1479
$ws.Range("N2").Value = 1479
$ws.Range("P2").Value = 44
$ws.Range("Q2").Value = 219
$ws.Range("R2").Value = -9
$ws.Range("S2").Value = -200
$ws.Range("T2").Value = 78
$ws.Range("U2").Value = 140
$ws.Range("V2").Value = 351
$ws.Range("W2").Value = 3.66
$ws.Range("X2").Value = 2.8
$ws.Range("Y2").Value = 4.82
$ws.Range("Z2").Value = 2.93
$ws.Range("AA2").Value = 58.99
$ws.Range("AB2").Value = 3871.17
$ws.Range("AC2").Value = 1577
$ws.Range("AD2").Value = 21.57
$ws.Range("AE2").Value = 47817
$ws.Range("AF2").Value = 0.71
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 0.88
$ws.Range("AI2").Value = 13.25
$ws.Range("AJ2").Value = 4444000

# Row 3: corrected financial figures
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = -7
$ws.Range("F3").Value = 123
$ws.Range("G3").Value = 51
$ws.Range("H3").Value = 140
$ws.Range("I3").Value = 140
$ws.Range("K3").Value = 2573
$ws.Range("L3").Value = 992
$ws.Range("M3").Value = 1581
$ws.Range("N3").Value = 1581
$ws.Range("P3").Value = 44
$ws.Range("Q3").Value = -29
$ws.Range("R3").Value = -205
$ws.Range("S3").Value = 219
$ws.Range("T3").Value = 77
$ws.Range("U3").Value = -105
$ws.Range("V3").Value = 580
$ws.Range("W3").Value = -191.6
$ws.Range("X3").Value = 3923.34
$ws.Range("Y3").Value = 9.130000000000001
$ws.Range("Z3").Value = 5.67
$ws.Range("AA3").Value = 62.76
$ws.Range("AB3").Value = 4164.23
$ws.Range("AC3").Value = 3143
$ws.Range("AD3").Value = 13.04
$ws.Range("AE3").Value = 51101
$ws.Range("AF3").Value = 0.8
$ws.Range("AG3").Value = 300
$ws.Range("AH3").Value = 0.73
$ws.Range("AI3").Value = 6.64
$ws.Range("AJ3").Value = 4444000

# Row 4: corrected financial figures
$ws.Range("D4").Value = 203
$ws.Range("E4").Value = 14
$ws.Range("F4").Value = 14
$ws.Range("G4").Value = -5
$ws.Range("H4").Value = 16
$ws.Range("I4").Value = 13
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 1186
$ws.Range("L4").Value = 102
$ws.Range("M4").Value = 1084
$ws.Range("N4").Value = 1070
$ws.Range("O4").Value = 14
$ws.Range("P4").Value = 22
$ws.Range("Q4").Value = 51
$ws.Range("R4").Value = -120
$ws.Range("S4").Value = 7
$ws.Range("T4").Value = 83
$ws.Range("U4").Value = -33
$ws.Range("V4").Value = 18
$ws.Range("W4").Value = 6.91
$ws.Range("X4").Value = 8.029999999999999
$ws.Range("Y4").Value = 0.97
$ws.Range("Z4").Value = 0.87
$ws.Range("AA4").Value = 9.380000000000001
$ws.Range("AB4").Value = 8584.530000000001
$ws.Range("AC4").Value = 387
$ws.Range("AD4").Value = 94.92
$ws.Range("AE4").Value = 71241
$ws.Range("AF4").Value = 0.52
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 0.54
$ws.Range("AI4").Value = 23.4
$ws.Range("AJ4").Value = 2159857

# Row 5: corrected financial figures
$ws.Range("D5").Value = 2634
$ws.Range("E5").Value = 180
$ws.Range("F5").Value = 180
$ws.Range("G5").Value = 195
$ws.Range("H5").Value = 128
$ws.Range("I5").Value = 79
$ws.Range("J5").Value = 49
$ws.Range("K5").Value = 3116
$ws.Range("L5").Value = 948
$ws.Range("M5").Value = 2168
$ws.Range("N5").Value = 1415
$ws.Range("O5").Value = 753
$ws.Range("P5").Value = 29
$ws.Range("Q5").Value = 328
$ws.Range("R5").Value = -436
$ws.Range("S5").Value = 116
$ws.Range("T5").Value = 129
$ws.Range("U5").Value = 199
$ws.Range("V5").Value = 380
$ws.Range("W5").Value = 6.83
$ws.Range("X5").Value = 4.87
$ws.Range("Y5").Value = 6.37
$ws.Range("Z5").Value = 5.96
$ws.Range("AA5").Value = 43.71
$ws.Range("AB5").Value = 7651.37
$ws.Range("AC5").Value = 2788
$ws.Range("AD5").Value = 10.83
$ws.Range("AE5").Value = 63803
$ws.Range("AF5").Value = 0.47
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 0.66
$ws.Range("AI5").Value = 5.6
$ws.Range("AJ5").Value = 2875800

# Row 6: corrected financial figures
$ws.Range("D6").Value = 2759
$ws.Range("E6").Value = 197
$ws.Range("F6").Value = 197
$ws.Range("G6").Value = 271
$ws.Range("H6").Value = 212
$ws.Range("I6").Value = 129
$ws.Range("K6").Value = 3209
$ws.Range("L6").Value = 879
$ws.Range("M6").Value = 2330
$ws.Range("N6").Value = 1500
$ws.Range("P6").Value = 29
$ws.Range("Q6").Value = 276
$ws.Range("R6").Value = -171
$ws.Range("S6").Value = -89
$ws.Range("T6").Value = 177
$ws.Range("U6").Value = 99
$ws.Range("V6").Value = 306
$ws.Range("W6").Value = 7.13
$ws.Range("X6").Value = 7.7
$ws.Range("Y6").Value = 8.82
$ws.Range("Z6").Value = 6.72
$ws.Range("AA6").Value = 37.71
$ws.Range("AB6").Value = 8077.82
$ws.Range("AC6").Value = 4470
$ws.Range("AD6").Value = 7.48
$ws.Range("AE6").Value = 67652
$ws.Range("AF6").Value = 0.49
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 0.6
$ws.Range("AI6").Value = 3.45
$ws.Range("AJ6").Value = 2875800

# Rows 7-9 (2019E/2020E/2021E estimate rows): remove all projected figures,
# keeping only the row index / period / company-name columns (A:C).
$ws.Range("D7:AI9").ClearContents()
